$d = $word.ActiveDocument

# --- 1. Remove the "_GoBack" bookmark currently sitting after "...,   GIT" ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- 2. Remove "AngularJS" (and its trailing ",   " separator) from the
#        INTERESTS line: "...Soccer,   AngularJS,   Kotlin,   ..."
#        becomes                "...Soccer,   Kotlin,   ..."
$find1 = $d.Content
[void]$find1.Find.Execute("AngularJS", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$angularStart = $find1.Start
$angularEnd = $find1.End

$angularRange = $d.Range($angularStart, $angularEnd)
$angularRange.Delete()

$find2 = $d.Content
[void]$find2.Find.Execute("Kotlin", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$kotlinStart = $find2.Start

$sepRange = $d.Range($angularStart, $kotlinStart)
$sepRange.Delete()

# --- 3. Word relocates the "_GoBack" bookmark to mark the spot of the last
#        edit, i.e. right before "Kotlin" (where "AngularJS" used to be). ---
$find3 = $d.Content
[void]$find3.Find.Execute("Kotlin", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$newBookmarkRange = $d.Range($find3.Start, $find3.Start)
$d.Bookmarks.Add("_GoBack", $newBookmarkRange)
